# Update the "Subject" column (B) with more specific subject names,
# replacing the old generic "Mathematics"/"Language"/"Elective" labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "Math"
    3  = "English"
    7  = "Fine Arts"
    8  = "English"
    11 = "English"
    13 = "English"
    15 = "Fine Arts"
    16 = "Math"
    19 = "Fine Arts"
    20 = "Personal Finance"
    23 = "Practical Arts"
    24 = "Math"
    25 = "Math"
    26 = "English"
    30 = "Math"
    31 = "English"
    32 = "Fine Arts"
    35 = "Math"
    37 = "Math"
    38 = "Math"
    39 = "English"
    43 = "Fine Arts"
    44 = "English"
    46 = "Math"
    47 = "Math"
    49 = "English"
    50 = "Practical Arts"
    53 = "Math"
    57 = "English"
    58 = "English"
    62 = "Practical Arts"
    63 = "Math"
    65 = "English"
    68 = "Math"
    71 = "Math"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# Reflect where the author's selection ended up after editing.
$ws.Range("B73").Select()
